# BRS-892-2 repairs to migration 2022/11/02
#
# Fixes a typo in the "Park Sub Area" label for Wallace Island Marine Park,
# corrects the Day Use / Boating flags and the Sub Area ID for Wallace
# Island Marine Park (row 2), corrects the Sub Area ID for Bodega Ridge
# Park (row 3), and normalizes row 2's cell formatting to match the rest
# of the sheet (Calibri 11 style already used by the other cells in the
# row), including the taller row height that comes with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections -------------------------------------------------

# Typo: "Conover Cover" -> "Conover Cove"
$ws.Range("G2").Value = "Wallace Island Marine Park - Conover Cove"

# Wallace Island Marine Park actually offers Day Use and Boating
$ws.Range("K2").Value = "Yes"
$ws.Range("L2").Value = "Yes"

# Sub Area IDs were mixed up between the two parks
$ws.Range("O2").Value = "0473"
$ws.Range("O3").Value = "0523"

# --- Formatting: bring row 2 in line with the rest of the sheet -------
# (reuse the existing "Calibri 11" cell style already applied to B2:E2/G2
# instead of inventing a new style entry)

$fmtSrc = $ws.Range("B2")
$fmtSrc.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$fmtSrc.Copy()
$ws.Range("F2").PasteSpecial(-4122)
$fmtSrc.Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(2).RowHeight = 15.75

# --- Cosmetic: cursor position -----------------------------------------
[void]$ws.Range("O8").Select()
